$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "Generation 0" row (row 8). Remaining rows below shift up,
# carrying their formatting with them (true row delete).
$ws.Rows.Item(8).Delete()

# After the shift, the rows that used to read "2020-2021" (Generation I,
# Generation I.5, Generation II) are now rows 9, 10 and 11 - update their
# Date column (E) to "2019-2021".
$ws.Range("E9").Value = "2019-2021"
$ws.Range("E10").Value = "2019-2021"
$ws.Range("E11").Value = "2019-2021"

# Update the selection / active cell as recorded in the saved workbook.
$ws.Range("F13").Select()

# Reposition the workbook window, matching the saved xWindow/yWindow.
try {
    $win = $excel.Windows.Item(1)
    $win.Left = 2730
    $win.Top = 2730
} catch {}
